# Fig3.xlsx – "2017-02-13 snapshot" update
# Rolls the workbook forward from the January 2017 STEO snapshot to the
# February 2017 snapshot: updates the report-title / source strings,
# refreshes the diesel/crude-oil price history + forecast-bracket cells,
# and lets Excel's own recalculation engine regenerate every dependent
# formula (column E "Difference") and chart cache.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fig3")

# ---------------------------------------------------------------------
# 1. Report title / source strings (January 2017 -> February 2017)
# ---------------------------------------------------------------------
$ws.Range("A2").Value   = "Short-Term Energy Outlook, February 2017"
$ws.Range("A100").Value = "Source: Short-Term Energy Outlook, February 2017."

# ---------------------------------------------------------------------
# 2. Updated history / forecast data, columns B (Diesel) and C (Cost of
#    Oil). Column E ("Difference") is the formula B-C and recalculates
#    automatically.
# ---------------------------------------------------------------------
$data = @(
    @{Row=72; C=1.0173809523809523},
    @{Row=73; C=1.1147619047619048},
    @{Row=74; C=1.0626190476190476},
    @{Row=76; B=2.5798000000000001; C=1.2259523809523809},
    @{Row=77; B=2.675252; C=1.2380952380952381},
    @{Row=78; B=2.7418040000000001; C=1.2380952380952381},
    @{Row=79; B=2.7054340000000003; C=1.2380952380952381},
    @{Row=80; B=2.6928390000000002; C=1.2380952380952381},
    @{Row=81; B=2.6931619999999996; C=1.2380952380952381},
    @{Row=82; B=2.7024209999999997; C=1.2619047619047619},
    @{Row=83; B=2.7296899999999997; C=1.2619047619047619},
    @{Row=84; B=2.758562; C=1.2619047619047619},
    @{Row=85; B=2.7683879999999998; C=1.2619047619047619},
    @{Row=86; B=2.786724; C=1.2619047619047619},
    @{Row=87; B=2.7932920000000001; C=1.2619047619047619},
    @{Row=88; B=2.7648440000000001; C=1.2619047619047619},
    @{Row=89; B=2.79555; C=1.2619047619047619},
    @{Row=90; B=2.8433679999999999; C=1.2619047619047619},
    @{Row=91; B=2.8068020000000002; C=1.2857142857142858},
    @{Row=92; B=2.8107889999999998; C=1.3095238095238095},
    @{Row=93; B=2.8202069999999999; C=1.3095238095238095},
    @{Row=94; B=2.8197260000000002; C=1.3095238095238095},
    @{Row=95; B=2.843483; C=1.3333333333333333},
    @{Row=96; B=2.8769429999999998; C=1.3333333333333333},
    @{Row=97; B=2.901586; C=1.3571428571428572},
    @{Row=98; B=2.9243700000000001; C=1.3571428571428572},
    @{Row=99; B=2.9455220000000004; C=1.3809523809523809}
)

foreach ($item in $data) {
    $r = $item.Row
    if ($item.ContainsKey("B")) {
        $ws.Cells.Item($r, 2).Value = $item.B
    }
    if ($item.ContainsKey("C")) {
        $ws.Cells.Item($r, 3).Value = $item.C
    }
}

# ---------------------------------------------------------------------
# 3. Forecast-bracket marker (scatter series) moves from month 48 to 49
# ---------------------------------------------------------------------
$ws.Range("A104").Value = 49
$ws.Range("A105").Value = 49

# ---------------------------------------------------------------------
# 4. Recalculate so column E and the chart caches pick up the new data
# ---------------------------------------------------------------------
$excel.CalculateFullRebuild()
